# Insert a new weekly data row at row 416, pushing existing rows 416-460 down to 417-461.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(416).Insert()

$ws.Range("A416").Value = 9
$ws.Range("B416").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C416").Value = "Metropolitana"
$ws.Range("D416").Value = 44769
$ws.Range("E416").Value = 13
$ws.Range("F416").Value = "Fruta"
$ws.Range("G416").Value = 100108
$ws.Range("H416").Value = "Tropicales y subtropicales"
$ws.Range("I416").Value = 100108002
$ws.Range("J416").Value = "Mango"
$ws.Range("K416").Value = "Sin especificar"
$ws.Range("L416").Value = "Especial"
$ws.Range("M416").Value = 200
$ws.Range("N416").Value = 9500
$ws.Range("O416").Value = 10000
$ws.Range("P416").Value = 9700
$ws.Range("Q416").Value = "`$/bandeja 4 kilos"
$ws.Range("R416").Value = "Brasil"
$ws.Range("S416").Value = 2425
$ws.Range("T416").Value = 4
